# Weekly driver report update for 2025-04-28
#
# The "Bad Drivers" table loses the "22.250.1.2" entry (a client count of 2
# moved below the >99.8% threshold it seems no longer qualifies as "bad" /
# was dropped), the remaining bad-driver rows are re-sorted and their
# figures refreshed, and several "Good Drivers" sample counts grow with a
# week's worth of new telemetry. Because one row disappears from the Bad
# Drivers block, everything beneath it (the blank spacer rows and the
# entire Good Drivers block) shifts up by a single row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the "Intel(R) Wi-Fi 6E AX211 160MHz - 22.250.1.2" row entirely.
#    Deleting the whole row shifts everything below it up by one, which is
#    exactly what the diff shows happening to the rest of the sheet.
$ws.Rows(4).Delete()

# 2. Rewrite the (now 3-row) Bad Drivers table with the refreshed,
#    re-sorted figures.
$ws.Cells.Item(3, 1).Value2 = "Intel(R) Wi-Fi 6E AX211 160MHz - 23.60.0.10"
$ws.Cells.Item(3, 2).Value2 = 2
$ws.Cells.Item(3, 3).Value2 = 249
$ws.Cells.Item(3, 4).Value2 = 97.09999999999999

$ws.Cells.Item(4, 1).Value2 = "Intel(R) Dual Band Wireless-AC 7265 - 19.51.48.1"
$ws.Cells.Item(4, 2).Value2 = 1
$ws.Cells.Item(4, 3).Value2 = 11
$ws.Cells.Item(4, 4).Value2 = 97.7

$ws.Cells.Item(5, 1).Value2 = "Intel(R) Wi-Fi 6E AX211 160MHz - 23.30.0.6"
$ws.Cells.Item(5, 2).Value2 = 164
$ws.Cells.Item(5, 3).Value2 = 16143
$ws.Cells.Item(5, 4).Value2 = 98.3

# 3. Update the Totals row (now row 6 after the row shift).
$ws.Cells.Item(6, 1).Value2 = "Totals:"
$ws.Cells.Item(6, 2).Value2 = 167
$ws.Cells.Item(6, 3).Value2 = 16403

# 4. Refresh the sample counts for five Good Drivers (their row numbers
#    shifted up by one thanks to the deletion above, but their order and
#    every other field stay the same).
$ws.Cells.Item(17, 2).Value2 = 11140    # 22.150.3.1: 10661 -> 11140
$ws.Cells.Item(18, 2).Value2 = 14487    # 22.150.0.3: 14239 -> 14487
$ws.Cells.Item(20, 2).Value2 = 236471   # 19.51.30.1: 205276 -> 236471
$ws.Cells.Item(21, 2).Value2 = 41618    # 19.51.29.1: 40211 -> 41618
$ws.Cells.Item(22, 2).Value2 = 132352   # 19.51.14.1: 122297 -> 132352
